$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells stay text even for values that look numeric
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "48.775.15"
$ws.Range("E2").Value = "  -1.68%  "
$ws.Range("D3").Value = "2.624.78"
$ws.Range("E3").Value = "  +3.15%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "109.48"
$ws.Range("E5").Value = "  +0.62%  "
$ws.Range("D6").Value = "322.03"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("E7").Value = "  -1.28%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "0.539"
$ws.Range("E9").Value = "  -3.23%  "
$ws.Range("D10").Value = "39.26"
$ws.Range("E10").Value = "  -2.14%  "
$ws.Range("D11").Value = "19.72"
$ws.Range("E11").Value = "  -2.94%  "
$ws.Range("D12").Value = "0.0808"
$ws.Range("E12").Value = "  -0.76%  "
$ws.Range("E13").Value = "  +0.13%  "
$ws.Range("D14").Value = "7.20"
$ws.Range("E14").Value = "  -0.55%  "
$ws.Range("D15").Value = "3.044.58"
$ws.Range("E15").Value = "  +3.33%  "
$ws.Range("D16").Value = "2.627.00"
$ws.Range("E16").Value = "  +5.19%  "
$ws.Range("D17").Value = "0.860"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("D18").Value = "48.746.99"
$ws.Range("E18").Value = "  -1.50%  "
$ws.Range("D19").Value = "12.75"
$ws.Range("E19").Value = "  -3.19%  "
$ws.Range("D20").Value = "2.92"
$ws.Range("E20").Value = "  -4.05%  "
$ws.Range("D21").Value = "6.67"
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("D23").Value = "269.21"
$ws.Range("E23").Value = "  -5.83%  "
$ws.Range("D24").Value = "68.58"
$ws.Range("E24").Value = "  -4.57%  "
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("D26").Value = "26.02"
$ws.Range("E26").Value = "  -1.29%  "
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").Value = "10.00"
$ws.Range("E28").Value = "  +1.92%  "
$ws.Range("E29").Value = "  -0.39%  "
$ws.Range("E30").Value = "  -4.42%  "
$ws.Range("D31").Value = "34.75"
$ws.Range("E31").Value = "  -1.10%  "
$ws.Range("D32").Value = "49.15"
$ws.Range("E32").Value = "  -0.79%  "
$ws.Range("D33").Value = "5.46"
$ws.Range("E33").Value = "  +1.77%  "
$ws.Range("D34").Value = "19.21"
$ws.Range("E34").Value = "  -1.83%  "
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("E36").Value = "  +0.94%  "
$ws.Range("E37").Value = "  +5.22%  "
$ws.Range("D38").Value = "2.03"
$ws.Range("E38").Value = "  +1.17%  "
$ws.Range("E39").Value = "  +4.69%  "
$ws.Range("D40").Value = "125.15"
$ws.Range("E40").Value = "  +3.38%  "
$ws.Range("D41").Value = "22.42"
$ws.Range("E41").Value = "  +0.65%  "
$ws.Range("E42").Value = "  -1.06%  "
$ws.Range("D43").Value = "2.15"
$ws.Range("E43").Value = "  -3.11%  "
$ws.Range("E44").Value = "  +0.89%  "
$ws.Range("D45").Value = "2.056.94"
$ws.Range("E45").Value = "  +2.28%  "
$ws.Range("D46").Value = "3.21"
$ws.Range("E46").Value = "  -1.95%  "
$ws.Range("E47").Value = "  +4.56%  "
$ws.Range("E48").Value = "  +1.79%  "
$ws.Range("D49").Value = "8.91"
$ws.Range("E49").Value = "  -1.14%  "
$ws.Range("D50").Value = "58.39"
$ws.Range("E50").Value = "  +2.27%  "
$ws.Range("D51").Value = "5.16"
$ws.Range("E51").Value = "  -2.77%  "
